$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 64, shifting existing rows 64-149 down to 65-150.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record.
$ws.Cells.Item(64, 1).Value = 5
$ws.Cells.Item(64, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(64, 3).Value = "Maule"
$ws.Cells.Item(64, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(64, 5).Value = 7
$ws.Cells.Item(64, 6).Value = "Fruta"
$ws.Cells.Item(64, 7).Value = 100108
$ws.Cells.Item(64, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(64, 9).Value = 100108005
$ws.Cells.Item(64, 10).Value = "Piña"
$ws.Cells.Item(64, 11).Value = "Caramelo"
$ws.Cells.Item(64, 12).Value = "Segunda"
$ws.Cells.Item(64, 13).Value = 150
$ws.Cells.Item(64, 14).Value = 19000
$ws.Cells.Item(64, 15).Value = 19000
$ws.Cells.Item(64, 16).Value = 19000
$ws.Cells.Item(64, 17).Value = '$/caja 14 unidades'
$ws.Cells.Item(64, 18).Value = "Ecuador"
$ws.Cells.Item(64, 19).Value = 1357
$ws.Cells.Item(64, 20).Value = 14
